$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览 = sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 7351
$wsExhibit.Range("F11").Value = 7512
$wsExhibit.Range("F14").Value = 5954
$wsExhibit.Range("F15").Value = 3211
$wsExhibit.Range("F20").Value = 429
$wsExhibit.Range("F22").Value = 263
$wsExhibit.Range("F23").Value = 2061
$wsExhibit.Range("F30").Value = 2563
$wsExhibit.Range("F31").Value = 1383
$wsExhibit.Range("F35").Value = 220

# Sheet "演出" (演出 = sheet2.xml)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F13").Value = 9

# Sheet "全部类型" (全部类型 = sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F15").Value = 7351
$wsAll.Range("F17").Value = 7512
$wsAll.Range("F20").Value = 5954
$wsAll.Range("F24").Value = 429
$wsAll.Range("F28").Value = 263
$wsAll.Range("F29").Value = 2061
$wsAll.Range("F30").Value = 9
$wsAll.Range("F37").Value = 2563
$wsAll.Range("F38").Value = 1383
$wsAll.Range("F42").Value = 220
